$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.830.25"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "3.023.03"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.018.20"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "3.513.53"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "62.015.95"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").Value = "3.029.62"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.675"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "459.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.85%  "
$ws.Range("D38").Value = "3.204.65"
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0798"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0384"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.29%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("D50").Value = "0.0₃0497"
$ws.Range("E50").Value = "  -5.86%  "
$ws.Range("E51").Value = "  +7.17%  "
